# DMS: Translate LuckyNumber Template
# Rename worksheets from their English identifiers to the Vietnamese
# labels used throughout the template's own content (sharedStrings
# already contain "Giải thưởng", "Mã", "Tên", etc.).

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("LuckyNumber").Name = "Giải thưởng"
$wb.Worksheets.Item("Organization").Name = "Đơn vị tổ chức"
